$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "46.0 (36.3+6.36,-8.09)"
$ws.Range("C2").Value = "207.0 (115.0+41.1,-28.8)"
$ws.Range("D2").Value = "48.0 (83.5+19.3,-27.7)"
$ws.Range("E2").Value = 9.869999999999999
